$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts the existing rows 6-9 down to
# 7-10 (carrying their styles/merges along), bumps the used-range dimension
# to H10, and slides the A6:A7 merged cell down to A7:A8 automatically.
$ws.Rows.Item(6).Insert()

# The newly inserted row 6 cloned its formatting from row 5 above it, so
# A6/C6/H6 need their border styles fixed up to match the "inner" row look
# (style 11 / 11 / 9) used by the rest of the data rows. Copy single cells
# that already carry those styles onto the blank ones; the values will be
# overwritten right after.
$ws.Range("A9").Copy($ws.Range("A6"))
$ws.Range("C9").Copy($ws.Range("C6"))
$ws.Range("H7").Copy($ws.Range("H6"))

# Row 5 (cyl = 4, vs = 1 group)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 91
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 2.14
$ws.Range("G5").ClearContents()

# Row 6 (cyl = 4, vs = 0 group)
$ws.Range("A6").Value = 4
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 81.8
$ws.Range("E6").Value = 21.87235698318771
$ws.Range("F6").Value = 2.3003
$ws.Range("G6").Value = 0.5982073312080948

# Row 7 (cyl = 6, vs = 0 group)
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 131.6666666666667
$ws.Range("E7").Value = 37.52776749732568
$ws.Range("F7").Value = 2.755
$ws.Range("G7").Value = 0.1281600561797629

# Row 8 (cyl = 6, vs = 1 group)
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 115.25
$ws.Range("E8").Value = 9.178779875342908
$ws.Range("F8").Value = 3.38875
$ws.Range("G8").Value = 0.1162163929916946

# Row 9 (cyl = 8 group)
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 209.2142857142857
$ws.Range("E9").Value = 50.97688551827051
$ws.Range("F9").Value = 3.999214285714287
$ws.Range("G9").Value = 0.7594047444769265

Write-Output "done"
